$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F3: empty -> 17.64
$ws.Range("F3").Value = 17.64

# Update F5: 17.66 -> empty
$ws.Range("F5").Value = $null

# Update F21: empty -> 16.58
$ws.Range("F21").Value = 16.58

# Update F23: 16.48 -> empty
$ws.Range("F23").Value = $null

# Delete entire row 26 (RM 232) - shifts rows up
$ws.Rows.Item(26).Delete()

# After deleting row 26, original row 28 (SC 92) is now row 27. Delete it.
$ws.Rows.Item(27).Delete()

# After both deletions, original row 34 (SC 193) is now row 32. Update its F value.
$ws.Range("F32").Value = 17.39
